$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT even when it looks like a plain number
# (e.g. "213.47"), by routing it through a scratch cell as a text-literal
# formula and pasting only the value. A direct `.Value = "213.47"` would be
# auto-coerced to a numeric cell, same as typing it into Excel would.
function Set-TextValue($cellRef, $textValue) {
    $ws.Range("ZZ1").Formula = "=""" + $textValue + """"
    $ws.Range("ZZ1").Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    $ws.Range("ZZ1").ClearContents()
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "28.270.77"
$ws.Range("E2").Value = "  +2.65%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.586.65"
$ws.Range("E3").Value = "  +1.52%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +1.15%  "

# Row 5 - BNB
Set-TextValue "D5" "213.47"
$ws.Range("E5").Value = "  +1.33%  "

# Row 6 - XRP
Set-TextValue "D6" "0.491"
$ws.Range("E6").Value = "  +0.08%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +1.15%  "

# Row 8 - Solana
Set-TextValue "D8" "24.02"
$ws.Range("E8").Value = "  +7.00%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  -0.06%  "

# Row 10 - Dogecoin
$ws.Range("E10").Value = "  +0.31%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0887"
$ws.Range("E11").Value = "  +2.27%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.813.42"
$ws.Range("E12").Value = "  +1.53%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.595.22"
$ws.Range("E13").Value = "  +2.16%  "

# Row 14 - Polygon
$ws.Range("E14").Value = "  +1.86%  "

# Row 15 - Polkadot
$ws.Range("E15").Value = "  -0.20%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "28.265.01"
$ws.Range("E16").Value = "  +2.71%  "

# Row 17 - Litecoin
Set-TextValue "D17" "63.18"
$ws.Range("E17").Value = "  +1.03%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "227.98"

# Row 19 - ShibaInu
$ws.Range("E19").Value = "  -0.05%  "

# Row 20
$ws.Range("E20").Value = "  -0.67%  "

# Row 21
$ws.Range("E21").Value = "  +1.13%  "

# Row 22
$ws.Range("E22").Value = "  -1.46%  "

# Row 23
$ws.Range("E23").Value = "  -0.80%  "

# Row 24
$ws.Range("E24").Value = "  +0.74%  "

# Row 25
Set-TextValue "D25" "151.92"
$ws.Range("E25").Value = "  +1.31%  "

# Row 26
$ws.Range("E26").Value = "  -0.04%  "

# Row 27 & 28 - Stellar/Cosmos swap
$ws.Range("B27").Value = "Cosmos"
$ws.Range("C27").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D27" "6.56"
$ws.Range("E27").Value = "  -0.89%  "

$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D28" "0.107"
$ws.Range("E28").Value = "  -1.22%  "

# Row 29
$ws.Range("E29").Value = "  +1.12%  "

# Row 30
$ws.Range("E30").Value = "  +0.04%  "

# Row 31
$ws.Range("E31").Value = "  +0.30%  "

# Row 32
$ws.Range("E32").Value = "  -0.05%  "

# Row 33
Set-TextValue "D33" "3.16"
$ws.Range("E33").Value = "  -0.50%  "

# Row 34
$ws.Range("D34").Value = "1.401.62"
$ws.Range("E34").Value = "  -3.83%  "

# Row 35
$ws.Range("E35").Value = "  -1.35%  "

# Row 36
Set-TextValue "D36" "1.01"
$ws.Range("E36").Value = "  -7.99%  "

# Row 37
$ws.Range("E37").Value = "  +2.06%  "

# Row 38
$ws.Range("E38").Value = "  -0.25%  "

# Row 39
$ws.Range("E39").Value = "  +8.99%  "

# Row 40
Set-TextValue "D40" "0.538"
$ws.Range("E40").Value = "  -0.53%  "

# Row 41
$ws.Range("E41").Value = "  -0.62%  "

# Row 42
$ws.Range("E42").Value = "  +1.08%  "

# Row 43
$ws.Range("E43").Value = "  +0.74%  "

# Row 44
Set-TextValue "D44" "5.61"
$ws.Range("E44").Value = "  -1.76%  "

# Row 45
Set-TextValue "D45" "0.980"
$ws.Range("E45").Value = "  +0.97%  "

# Row 46
Set-TextValue "D46" "64.29"
$ws.Range("E46").Value = "  -0.99%  "

# Row 47
$ws.Range("D47").Value = "1.723.36"
$ws.Range("E47").Value = "  +1.30%  "

# Row 48
$ws.Range("E48").Value = "  +1.48%  "

# Row 49
Set-TextValue "D49" "86.70"
$ws.Range("E49").Value = "  +0.29%  "

# Row 50
$ws.Range("E50").Value = "  +5.57%  "

# Row 51
$ws.Range("E51").Value = "  -0.53%  "
